# Apply the change described by the diff for "Översikt ORSA.xlsx":
#  1. Column C ("Förändrad") for every existing data row (2-288) changes
#     from 45190 to 45192.
#  2. Row 288 gains an explicit row height (ht="15" customHeight="1").
#  3. A brand-new data row (289) is appended with the new ORSA case
#     "A 44600-2023".
#  4. The sheet dimension grows from A1:Y288 to A1:Y289 (handled
#     automatically by Excel once the new row's cells are populated).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Update the "Förändrad" date for every existing row -----------------
$ws.Range("C2:C288").Value = 45192

# --- 2. Row 288 now carries an explicit row height --------------------------
$ws.Rows.Item(288).RowHeight = 15

# --- 3. Append the new row 289 ----------------------------------------------
$ws.Range("A289").Value = "A 44600-2023"

$ws.Range("B289").Value = 45189
$ws.Range("B289").NumberFormat = "YYYY-MM-DD"

$ws.Range("C289").Value = 45192
$ws.Range("C289").NumberFormat = "YYYY-MM-DD"

$ws.Range("D289").Value = "DALARNAS LÄN"
$ws.Range("E289").Value = "ORSA"

$ws.Range("G289").Value = 3.9
$ws.Range("H289").Value = 0
$ws.Range("I289").Value = 0
$ws.Range("J289").Value = 0
$ws.Range("K289").Value = 0
$ws.Range("L289").Value = 0
$ws.Range("M289").Value = 0
$ws.Range("N289").Value = 0
$ws.Range("O289").Value = 0
$ws.Range("P289").Value = 0
$ws.Range("Q289").Value = 0

# R289 stays empty but keeps the wrap-text styling used throughout column R
$ws.Range("R289").Value = ""
$ws.Range("R289").WrapText = $true
